$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.987.80'
Set-TextValue 'E2' '  -0.24%  '

Set-TextValue 'D3' '1.868.21'

Set-TextValue 'E4' '  +0.13%  '

Set-TextValue 'D5' '318.13'
Set-TextValue 'E5' '  -1.28%  '

Set-TextValue 'D6' '1.002'
Set-TextValue 'E6' '  +0.10%  '

Set-TextValue 'D7' '0.5074'
Set-TextValue 'E7' '  -1.87%  '

Set-TextValue 'D8' '0.3915'
Set-TextValue 'E8' '  -2.25%  '

Set-TextValue 'D9' '0.08173'
Set-TextValue 'E9' '  -3.71%  '

Set-TextValue 'D10' '41.94'
Set-TextValue 'E10' '  -2.51%  '

Set-TextValue 'D11' '1.088'
Set-TextValue 'E11' '  -3.28%  '

Set-TextValue 'D12' '22.50'
Set-TextValue 'E12' '  +5.76%  '

Set-TextValue 'D13' '1.867.41'
Set-TextValue 'E13' '  -3.01%  '

Set-TextValue 'D14' '6.251'
Set-TextValue 'E14' '  -1.29%  '

Set-TextValue 'D15' '7.158'
Set-TextValue 'E15' '  -2.95%  '

Set-TextValue 'D16' '1.003'
Set-TextValue 'E16' '  +0.14%  '

Set-TextValue 'D17' '91.67'
Set-TextValue 'E17' '  -3.14%  '

Set-TextValue 'D18' '0.00001078'
Set-TextValue 'E18' '  -3.40%  '

Set-TextValue 'D19' '0.06321'
Set-TextValue 'E19' '  -6.40%  '

Set-TextValue 'D20' '17.81'
Set-TextValue 'E20' '  -1.06%  '

Set-TextValue 'D21' '1.003'
Set-TextValue 'E21' '  +0.23%  '

Set-TextValue 'D22' '29.979.62'
Set-TextValue 'E22' '  -0.34%  '

Set-TextValue 'D23' '5.778'
Set-TextValue 'E23' '  -4.87%  '

Set-TextValue 'D24' '11.04'
Set-TextValue 'E24' '  -1.57%  '

Set-TextValue 'D25' '2.205'
Set-TextValue 'E25' '  -0.13%  '

Set-TextValue 'D26' '2.087.75'
Set-TextValue 'E26' '  -2.66%  '

Set-TextValue 'D27' '161.12'
Set-TextValue 'E27' '  +0.65%  '

Set-TextValue 'E28' '  -0.98%  '

Set-TextValue 'D29' '2.223'
Set-TextValue 'E29' '  -10.04%  '

Set-TextValue 'D30' '126.63'
Set-TextValue 'E30' '  -2.25%  '

Set-TextValue 'D31' '0.1031'
Set-TextValue 'E31' '  -2.61%  '

Set-TextValue 'E32' '  -3.89%  '

Set-TextValue 'D33' '5.869'
Set-TextValue 'E33' '  -3.94%  '

Set-TextValue 'D34' '3.741'
Set-TextValue 'E34' '  +2.24%  '

Set-TextValue 'D35' '0.02413'
Set-TextValue 'E35' '  -3.46%  '

Set-TextValue 'D36' '5.191'
Set-TextValue 'E36' '  -0.22%  '

Set-TextValue 'D37' '0.06326'
Set-TextValue 'E37' '  -4.49%  '

Set-TextValue 'D38' '0.2129'
Set-TextValue 'E38' '  -4.09%  '

Set-TextValue 'B39' 'ARBITRUM'
Set-TextValue 'C39' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D39' '1.165'
Set-TextValue 'E39' '  -6.47%  '

Set-TextValue 'B40' 'FraxShare'
Set-TextValue 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '8.493'
Set-TextValue 'E40' '  -5.87%  '

Set-TextValue 'B41' 'TrustWalletToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D41' '1.211'
Set-TextValue 'E41' '  -2.33%  '

Set-TextValue 'B42' 'TheSandbox'
Set-TextValue 'C42' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D42' '0.6252'
Set-TextValue 'E42' '  -4.40%  '

Set-TextValue 'D43' '11.17'
Set-TextValue 'E43' '  -2.12%  '

Set-TextValue 'E44' '  +0.10%  '

Set-TextValue 'D45' '0.5869'
Set-TextValue 'E45' '  -4.45%  '

Set-TextValue 'D46' '12.78'
Set-TextValue 'E46' '  -3.27%  '

Set-TextValue 'D47' '3.633'
Set-TextValue 'E47' '  -2.37%  '

Set-TextValue 'D48' '1.984'
Set-TextValue 'E48' '  -3.67%  '

Set-TextValue 'D49' '121.80'
Set-TextValue 'E49' '  -3.01%  '

Set-TextValue 'D50' '1.201'
Set-TextValue 'E50' '  -3.36%  '

Set-TextValue 'E51' '  -2.56%  '
